# Trade #85 closed at 2026-02-17 15:53:24 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up numbers for the
# MarketMaking strategy and appends the newly-closed trade #85 row to
# both the "All Trades" and "MarketMaking" trade logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet roll-up figures
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.01   # Current Capital
$summary.Range("B4").Value = 0         # Total P&L $
$summary.Range("B5").Value = 0         # Total P&L %
$summary.Range("B6").Value = 85        # Total Trades
$summary.Range("B8").Value = 44        # Losing Trades
$summary.Range("B9").Value = 34.12     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.01
$status.Range("D4").Value = 85
$status.Range("E4").Value = 0
$status.Range("F4").Value = 0.01
$status.Range("G4").Value = 34.12

# ---------------------------------------------------------------------
# 3. Append trade #85 to both trade-log sheets (All Trades & MarketMaking)
# ---------------------------------------------------------------------
$tradeSheetNames = @("All Trades", "MarketMaking")

foreach ($name in $tradeSheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $row = 86

    $ws.Cells.Item($row, 1).Value = 85

    # Date/time columns are plain text in this workbook (not real Excel
    # dates) - force Text formatting on the date cell so the literal
    # string "2026-02-17" is preserved instead of being parsed into a
    # date serial number.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 3).Value = "15:53:17"

    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.24
    $ws.Cells.Item($row, 7).Value = 0.18
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -25
    $ws.Cells.Item($row, 10).Value = -0.06
    $ws.Cells.Item($row, 11).Value = 100.01
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.14
}
